$wb = $excel.ActiveWorkbook

# Sheet "CreateModifyDeleteProfile": latest Demo run timestamp for the single test row
$wsProfile = $wb.Worksheets.Item("CreateModifyDeleteProfile")
$wsProfile.Range("B2").Value = "Tue Aug 12 01:53:59 IST 2025"

# Sheet "AddModifyDeleteCC": latest Demo run timestamp for the single test row
$wsCC = $wb.Worksheets.Item("AddModifyDeleteCC")
$wsCC.Range("B2").Value = "Tue Aug 12 01:52:58 IST 2025"

# Sheet "AddModifyDeleteACH": latest Demo run timestamps for the three test rows
$wsACH = $wb.Worksheets.Item("AddModifyDeleteACH")
$wsACH.Range("B2").Value = "Tue Aug 12 01:49:40 IST 2025"
$wsACH.Range("B3").Value = "Tue Aug 12 01:50:50 IST 2025"
$wsACH.Range("B4").Value = "Tue Aug 12 01:51:54 IST 2025"
